# "Solved some more DP questions"
# Adds two new LeetCode entries (rows 45 and 46) below the existing table,
# removes the stray numeric value that was left in A44, and updates the
# sheet's dimension/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44: the plain number in A44 was a leftover artifact; clear it so
# the row only keeps its B (category text) and D (hyperlink) cells. ---
$null = $ws.Range("A44").ClearContents()

# --- Row 45: "Combination Sum IV" ---
# Set D (link) first, then B (title), then A (category) so that new shared
# strings are appended to sharedStrings.xml in the same order as the target.
$ws.Range("D45").Value2 = "https://leetcode.com/problems/combination-sum-iv/"
$ws.Range("B45").Value2 = "Combination Sum IV"
$ws.Range("A45").Value2 = "45. Array + DP"

# --- Row 46: "Search Suggestions System (June Day-19)" (no category cell) ---
$ws.Range("D46").Value2 = "https://leetcode.com/problems/search-suggestions-system/"
$ws.Range("B46").Value2 = "Search Suggestions System (June Day-19)"

# Match the "Good"-style formatting (with the sheet's custom fill override)
# used by the rest of column B, copying it from an existing, equivalently
# styled cell.
$null = $ws.Range("B40").Copy()
$null = $ws.Range("B45:B46").PasteSpecial(-4122)

# Create the actual hyperlinks for the new URLs.
$null = $ws.Hyperlinks.Add($ws.Range("D45"), "https://leetcode.com/problems/combination-sum-iv/")
$null = $ws.Hyperlinks.Add($ws.Range("D46"), "https://leetcode.com/problems/search-suggestions-system/")

# Adding a hyperlink re-applies its own formatting; restore column D to the
# exact same "Hyperlink" cell style already used throughout the sheet.
$hlStyle = $ws.Range("D44").Style()
$null = $ws.Range("D45").Style = $hlStyle
$null = $ws.Range("D46").Style = $hlStyle

# The new last populated cell is C48, matching the updated selection.
$null = $ws.Range("C48").Select()
